$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the new range as Text so the numeric-looking IDs are
# stored as shared strings (matching the source data) instead of numbers.
$ws.Range("A46:D54").NumberFormat = "@"

$ws.Range("A46").Value = "118451"
$ws.Range("B46").Value = "1008617900"
$ws.Range("C46").Value = "17706587"
$ws.Range("D46").Value = "1001"

$ws.Range("A47").Value = "118451"
$ws.Range("B47").Value = "1008617901"
$ws.Range("C47").Value = "17706588"
$ws.Range("D47").Value = "1001"

$ws.Range("A48").Value = "118452"
$ws.Range("B48").Value = "1008617904"
$ws.Range("C48").Value = "17706592"
$ws.Range("D48").Value = "1150"

$ws.Range("A49").Value = "118448"
$ws.Range("B49").Value = "1008617917"
$ws.Range("C49").Value = "17707507"
$ws.Range("D49").Value = "1010"

$ws.Range("A50").Value = "118448"
$ws.Range("B50").Value = "1008617918"
$ws.Range("C50").Value = "17707512"
$ws.Range("D50").Value = "1010"

$ws.Range("A51").Value = "118448"
$ws.Range("B51").Value = "1008617924"
$ws.Range("C51").Value = "17707515"
$ws.Range("D51").Value = "1010"

$ws.Range("A52").Value = "118451"
$ws.Range("B52").Value = "1008617926"
$ws.Range("C52").Value = "17707517"
$ws.Range("D52").Value = "1007"

$ws.Range("A53").Value = "118451"
$ws.Range("B53").Value = "1008617927"
$ws.Range("C53").Value = "17707519"
$ws.Range("D53").Value = "1010"

$ws.Range("A54").Value = "118452"
$ws.Range("B54").Value = "1008617928"
$ws.Range("C54").Value = "17707520"
$ws.Range("D54").Value = "1011"

# Remove the temporary Text number-format again so the cells are left
# with the default style, matching the rest of the sheet.
$ws.Range("A46:D54").ClearFormats()
